$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Park a copy of the existing file_id cell formatting (Tahoma 10) somewhere
# safe before the table below gets cleared out.
$ws.Range("D2").Copy($ws.Range("Z100"))

# The sticker table used to live in B1:E3. Clear it out completely (values
# and formatting -- this keeps the sheet's column-width definitions intact)
# and re-enter the data shifted one column to the left, into A1:D3, plus two
# new rows of sticker data appended below it.
$ws.Range("B1:E3").Clear()

# Header row.
$ws.Range("A1").Value = "emoji"
$ws.Range("B1").Value = "file_id"
$ws.Range("C1").Value = "sticker_id"
$ws.Range("D1").Value = "key_word"

# Existing two stickers, now in columns C (file_id) / D (key_word).
$ws.Range("C2").Value = "CAACAgIAAxkBAANBYFnKfcWD9t6m_8-4LD8clr4e4wcAAlkAAwr8wgXOU7sZfH5zGx4E"
$ws.Range("D2").Value = "смешная шутка"

$ws.Range("C3").Value = "CAACAgIAAxkBAANEYFnLQxjtNi5MTuMghLPi9mJjD3MAAg0EAALPX4sHmuYS8a7yxGQeBA"
$ws.Range("D3").Value = "я не хочу брать Иерусалим"

# Two new stickers appended as rows 4 and 5 (key_word typed before file_id
# for each row).
$ws.Range("D4").Value = "пока"
$ws.Range("C4").Value = "CAACAgIAAxkBAANKYFsYAAE4EvZcktpJ37Vholo6BLUVAAKKAgACVp29Cj5SbosTxUBnHgQ"

$ws.Range("D5").Value = "привет"
$ws.Range("C5").Value = "CAACAgIAAxkBAAOFYFse2EbZxvImJ_jaCeqYhMXZzUUAAkMAA4wPBgUBj55LMpPjEB4E"

# Give the whole file_id column (C2:C5) the established look (Tahoma 10,
# black) by pasting back the formatting we parked earlier.
$ws.Range("Z100").Copy()
$ws.Range("C2:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z100").Clear()

# Resize columns C and D to fit the new content (column E's old width is
# left as-is since there's no data there anymore).
$ws.Columns("C").ColumnWidth = 18.416666667
$ws.Columns("D").ColumnWidth = 24.583333334

# Leave the selection where the user finished editing.
$ws.Range("C5").Select()
